# jobs module: added subpackage for results
#
# Adds a "parent" column to the packages sheet and a new "jobs_results"
# sub-package row (child of "jobs"), describing a new results entity
# subpackage for the jobs module.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("packages")

# New header column D: "parent"
$ws.Cells.Item(1, 4).Value2 = "parent"

# New row describing the "jobs_results" package (child of "jobs")
$ws.Cells.Item(3, 1).Value2 = "jobs_results"
$ws.Cells.Item(3, 2).Value2 = "Results"
$ws.Cells.Item(3, 3).Value2 = "The outcome of process including analyzed data."
$ws.Cells.Item(3, 4).Value2 = "jobs"

Write-Output "packages sheet updated with parent column and jobs_results row"
